# Add a "Save" column (H) to the s_vals sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1: same text style as the other header cells (B1:G1).
$ws.Range("H1").Value = "Save"

# Copy the header formatting from the neighboring header cell (G1) onto H1
# so it reuses the existing header style instead of creating a new one.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# New data cell H2: plain numeric value, no special formatting (like B2:G2).
$ws.Range("H2").Value = 1
